$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.271.22"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "2.587.43"
$ws.Range("E3").Value = "  -3.91%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.86%  "

$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("D9").Value = "2.593.95"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.32%  "

$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").Value = "3.043.15"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "60.270.02"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.33%  "

$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "2.592.46"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "

$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "

$ws.Range("D28").Value = "0.0₃0840"
$ws.Range("E28").Value = "  -3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.857"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "299.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.03%  "

$ws.Range("E49").Value = "  -2.26%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "1.993.01"
$ws.Range("E51").Value = "  -2.29%  "
